$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# New header for column D
$ws.Range("D4").Value = "coda LE"

# New values for column D rows 5-10
$ws.Range("D5").Value = 81920
$ws.Range("D6").Value = 65536
$ws.Range("D7").Value = 622591
$ws.Range("D8").Value = 671743
$ws.Range("D9").Value = 589823
$ws.Range("D10").Value = 557055

# Update selection to match diff (F16)
$ws.Range("F16").Select()
